$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122 (shifts old rows 122..157 down to 123..158,
# and grows the used range from A1:R157 to A1:R158).
$ws.Rows("122:122").Insert()

# Populate the new row 122 with a fresh weekly record (same base attributes
# as the old row 122, which is now row 123, but with its own Fecha/Volumen).
$ws.Range("A122").Value = 10
$ws.Range("B122").Value = "Vega Modelo de Temuco"
$ws.Range("C122").Value = "La Araucanía"
$ws.Range("D122").Value = 44932
$ws.Range("E122").Value = 9
$ws.Range("F122").Value = 100114002
$ws.Range("G122").Value = "Camote"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 20
$ws.Range("K122").Value = 24000
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = 24000
$ws.Range("N122").Value = "$/malla 20 kilos"
$ws.Range("O122").Value = "Perú"
$ws.Range("P122").Value = 1200
$ws.Range("Q122").Value = 20
$ws.Range("R122").Value = "Hortaliza"
